# refs #882 Video Wall -> HSR Videowall
#
# The header row of the "Aenderungsgeschichte" sheet needed a touch more
# vertical room once the title text was updated, so the row height for
# row 1 was bumped up (and pinned as an explicit/custom height) to fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1).RowHeight = 24.75
